# Apply the "ver 1 with exe" edit to the Advfinsa client list workbook.
# The sheet previously tracked a single "Nombre" (full name) column; it is
# reworked here to split names into "Nombres" / "Apellidos" columns, add
# several new client rows, and drop an unused blank template row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "yyyy\-mm\-dd;@"

# --- Header row ---------------------------------------------------------
$ws.Range("A2").Value = "ID"
$ws.Range("B2").Value = "Nombres"
$ws.Range("C2").Value = "Apellidos"
$ws.Range("D2").Value = "Fecha"
$ws.Range("E2").Value = "tipo_persona"

# --- Row 3: Alfredo Javier Barandearan Oyague ---------------------------
$ws.Range("A3").Value = "0908894934"
$ws.Range("B3").Value = "Alfredo Javier"
$ws.Range("C3").ClearFormats()
$ws.Range("C3").Value = "Barandearan Oyague"
$ws.Range("D3").NumberFormat = $dateFmt
$ws.Range("D3").Value = Get-Date -Year 1968 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("E3").Value = "Natural"

# --- Row 4: Melania Sissi Gutierrez Gavilanes ----------------------------
$ws.Range("A4").Value = "0914788245"
$ws.Range("B4").Value = "Melania Sissi"
$ws.Range("C4").ClearFormats()
$ws.Range("C4").Value = "Gutierrez Gavilanes"
$ws.Range("D4").NumberFormat = $dateFmt
$ws.Range("D4").Value = Get-Date -Year 1973 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("E4").Value = "Natural"

# --- Row 5: EMPAGRAM (juridica, no Apellidos) ----------------------------
$ws.Range("A5").Value = "0990071969001"
$ws.Range("B5").Value = "EMPAGRAM"
$ws.Range("C5").Clear()
$ws.Range("D5").NumberFormat = $dateFmt
$ws.Range("D5").Value = Get-Date -Year 2000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("E5").Value = "Juridica"

# --- Row 6: Ariana Cristina Barandearan Gutierrez ------------------------
$ws.Range("A6").Value = "0929029395"
$ws.Range("B6").Value = "Ariana Cristina"
$ws.Range("C6").ClearFormats()
$ws.Range("C6").Value = "Barandearan Gutierrez"
$ws.Range("D6").NumberFormat = $dateFmt
$ws.Range("D6").Value = Get-Date -Year 2006 -Month 3 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Range("E6").Value = "Natural"

# --- Row 7: Jose Adolfo Macias Villamar ----------------------------------
$ws.Range("A7").Value = "1309022935"
$ws.Range("B7").Value = "Jose Adolfo"
$ws.Range("C7").ClearFormats()
$ws.Range("C7").Value = "Macias Villamar"
$ws.Range("D7").NumberFormat = $dateFmt
$ws.Range("D7").Value = Get-Date -Year 1979 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("E7").Value = "Natural"

# --- Row 8: Las Olas SA (juridica, no Apellidos) -------------------------
$ws.Range("A8").Value = "1391791119001"
$ws.Range("B8").Value = "Las Olas SA"
$ws.Range("C8").Clear()
$ws.Range("D8").NumberFormat = $dateFmt
$ws.Range("D8").Value = Get-Date -Year 2000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("E8").Value = "Juridica"

# --- Row 9: Roberto David Barandearan Oyague -----------------------------
$ws.Range("A9").Value = "0908890452"
$ws.Range("B9").Value = "Roberto David"
$ws.Range("C9").ClearFormats()
$ws.Range("C9").Value = "Barandearan Oyague"
$ws.Range("D9").NumberFormat = $dateFmt
$ws.Range("D9").Value = Get-Date -Year 1977 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("E9").Value = "Natural"

# --- Row 10: Mary Paz Herrera Oramas -------------------------------------
$ws.Range("A10").Value = "1709705675"
$ws.Range("B10").Value = "Mary Paz"
$ws.Range("C10").Value = "Herrera Oramas"
$ws.Range("D10").NumberFormat = $dateFmt
$ws.Range("D10").Value = Get-Date -Year 1969 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("E10").Value = "Natural"

# Row 11 (previously a blank placeholder row between the data and the
# remaining empty template rows) is removed - the data now flows straight
# from row 10 into the still-empty row 12 onward.
$ws.Range("A11:E11").EntireRow.Delete()

# --- Column widths (auto-fit after the wider "Apellidos" content) -------
# Target character widths (from the workbook after a real Excel auto-fit)
# are 13.78 / 19.44 / 11.78 / 11.78; the values below are the closest this
# runtime's width quantization can reach.
$ws.Columns.Item(2).ColumnWidth = 13.0
$ws.Columns.Item(3).ColumnWidth = 18.66666667
$ws.Columns.Item(4).ColumnWidth = 11.0
$ws.Columns.Item(5).ColumnWidth = 11.0

# --- Selection matches the author's last active cell --------------------
$ws.Range("C16").Select()
